# Auto-generated edit script: update crypto price/volume cells per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.806.95"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "3.126.35"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'532.72"
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("D6").Value = "'138.80"
$ws.Range("E6").Value = "  +0.78%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.125.40"
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("D9").Value = "'0.469"
$ws.Range("E9").Value = "  +4.78%  "
$ws.Range("D10").Value = "'7.31"
$ws.Range("E10").Value = "  +1.34%  "
$ws.Range("D11").Value = "'0.108"
$ws.Range("E11").Value = "  +0.36%  "
$ws.Range("D12").Value = "'0.415"
$ws.Range("E12").Value = "  +4.44%  "
$ws.Range("D13").Value = "3.660.27"
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("E14").Value = "  +1.34%  "
$ws.Range("D15").Value = "'25.54"
$ws.Range("E15").Value = "  +0.67%  "
$ws.Range("D16").Value = "'0.0000165"
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("D17").Value = "57.994.02"
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("D18").Value = "3.128.11"
$ws.Range("D19").Value = "'6.03"
$ws.Range("E19").Value = "  +1.19%  "
$ws.Range("D20").Value = "'12.72"
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("E21").Value = "  +2.40%  "
$ws.Range("D22").Value = "'361.21"
$ws.Range("E22").Value = "  +3.54%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").Value = "'69.06"
$ws.Range("E24").Value = "  +1.00%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("E26").Value = "  -0.51%  "
$ws.Range("E27").Value = "  +0.43%  "
$ws.Range("D28").Value = "0.0₃0876"
$ws.Range("E28").Value = "  -4.09%  "
$ws.Range("D29").Value = "'7.31"
$ws.Range("E29").Value = "  -1.81%  "
$ws.Range("E30").Value = "  -0.27%  "
$ws.Range("D31").Value = "'6.09"
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("E32").Value = "  +1.28%  "
$ws.Range("D33").Value = "'5.14"
$ws.Range("E33").Value = "  +3.40%  "
$ws.Range("E34").Value = "  -2.47%  "
$ws.Range("D35").Value = "'158.26"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  -1.48%  "
$ws.Range("D37").Value = "'25.84"
$ws.Range("E37").Value = "  -1.12%  "
$ws.Range("E38").Value = "  +2.02%  "
$ws.Range("E39").Value = "  +2.68%  "
$ws.Range("D40").Value = "'0.0674"
$ws.Range("E40").Value = "  +1.20%  "
$ws.Range("D41").Value = "2.489.19"
$ws.Range("E41").Value = "  +6.11%  "
$ws.Range("D42").Value = "'0.699"
$ws.Range("E42").Value = "  -0.29%  "
$ws.Range("D43").Value = "'4.00"
$ws.Range("E43").Value = "  -4.78%  "
$ws.Range("D44").Value = "'37.70"
$ws.Range("E44").Value = "  +3.32%  "
$ws.Range("D45").Value = "3.167.07"
$ws.Range("E45").Value = "  +0.35%  "
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").Value = "'0.0269"
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("D48").Value = "'0.988"
$ws.Range("E48").Value = "  +3.00%  "
$ws.Range("D49").Value = "'6.08"
$ws.Range("E49").Value = "  +0.65%  "
$ws.Range("D50").Value = "'19.80"
$ws.Range("E50").Value = "  -1.69%  "
$ws.Range("D51").Value = "'0.743"
$ws.Range("E51").Value = "  -2.62%  "
